# Update latest output (run 133)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1): rows 2-5 ---
$schedule = $wb.Worksheets.Item("Schedule")

# Row 2
$schedule.Range("B2").Value = 46042.22916666666
$schedule.Range("C2").Value = 5.5
$schedule.Range("D2").Value = 20.79
$schedule.Range("E2").Value = 626.1786959999998
$schedule.Range("F2").Value = 30.11922539682539

# Row 3
$schedule.Range("A3").Value = 46042.27083333334
$schedule.Range("C3").Value = 9.5
$schedule.Range("D3").Value = 35.91
$schedule.Range("E3").Value = 32.91044249999999
$schedule.Range("F3").Value = 0.9164701336675019

# Row 4
$schedule.Range("A4").Value = 46042.95833333334
$schedule.Range("C4").Value = 4
$schedule.Range("D4").Value = 15.12
$schedule.Range("E4").Value = 492.34987425
$schedule.Range("F4").Value = 32.56282237103174

# Row 5
$schedule.Range("E5").Value = -197.16784425
$schedule.Range("F5").Value = -5.795645039682539

# --- Sheet "Detailed" (sheet2): rows 12, 15, 43-97 ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("E12").Value = "ON"
$detailed.Range("E15").Value = "ON"

$detailed.Range("B43").Value = 56.98
$detailed.Range("B44").Value = 56.98

$detailed.Range("B45").Value = 80.64212000000001
$detailed.Range("C45").Value = "historical"

$detailed.Range("B46").Value = 77.94
$detailed.Range("C46").Value = "historical"
$detailed.Range("E46").Value = "OFF"

$detailed.Range("B47").Value = 78
$detailed.Range("E47").Value = "OFF"

$detailed.Range("B48").Value = 65.85371000000001
$detailed.Range("B49").Value = 65
$detailed.Range("B50").Value = 65
$detailed.Range("B51").Value = 65
$detailed.Range("B52").Value = 65
$detailed.Range("B53").Value = 65
$detailed.Range("B55").Value = 57.06049
$detailed.Range("B56").Value = 65
$detailed.Range("B57").Value = 65.46223000000001
$detailed.Range("B58").Value = 66.07064
$detailed.Range("B59").Value = 67.19774
$detailed.Range("B60").Value = 68.17549
$detailed.Range("B61").Value = 78
$detailed.Range("B62").Value = 83.83537
$detailed.Range("B63").Value = 65
$detailed.Range("B64").Value = 30.90576
$detailed.Range("B65").Value = 0.64885
$detailed.Range("B66").Value = -3.42354
$detailed.Range("B68").Value = -6.82305
$detailed.Range("B69").Value = -6.93715
$detailed.Range("B70").Value = -9.565659999999999
$detailed.Range("B71").Value = -9.700060000000001
$detailed.Range("B73").Value = -14
$detailed.Range("B74").Value = -14.68544
$detailed.Range("B75").Value = -22.78962
$detailed.Range("B76").Value = -19.39014
$detailed.Range("B78").Value = -20.82009
$detailed.Range("B80").Value = -23.5
$detailed.Range("B81").Value = -15.14329
$detailed.Range("B82").Value = -6.88086
$detailed.Range("B83").Value = -5.46656
$detailed.Range("B85").Value = 47.21538
$detailed.Range("B86").Value = 54.41338
$detailed.Range("B87").Value = 57.31
$detailed.Range("B88").Value = 79.95
$detailed.Range("B89").Value = 83.63907
$detailed.Range("B91").Value = 71.40000000000001
$detailed.Range("B93").Value = 71.0382
$detailed.Range("B94").Value = 57.09
$detailed.Range("B95").Value = 63.95073
$detailed.Range("B97").Value = 64.88197
